$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.2174858043788258
$ws.Range("AB2").Value = -979.3960140216795
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -979.3960140216795

# Row 3
$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1933261684273758
$ws.Range("AB3").Value = -753.5748113152024
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -753.5748113152024

# Row 4
$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.1756095954340679
$ws.Range("AB4").Value = -704.9156816189599
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -704.9156816189599

# Row 5
$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.2022216384210701
$ws.Range("AB5").Value = -970.9709717395569
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -970.9709717395569

# Row 6
$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.1989786663229278
$ws.Range("AB6").Value = -1238.089479342662
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1238.089479342662

# Row 7
$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1278335385115142
$ws.Range("AB7").Value = -568.7486896660039
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -568.7486896660039

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").ClearContents()

# Row 9
$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.1843783617471505
$ws.Range("AB9").Value = 830.3044563426573
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 830.3044563426573

# Row 10
$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = 0.1277118059246617
$ws.Range("AB10").Value = 497.8136215871573
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 497.8136215871573

# Row 11
$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.2537996662423591
$ws.Range("AB11").Value = 1018.778981192217
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 1018.778981192217

# Row 12
$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.2813298779931148
$ws.Range("AB12").Value = 1350.810660754116
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1350.810660754116

# Row 13
$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.2192293694956227
$ws.Range("AB13").Value = 1364.09385463943
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1364.09385463943

# Row 14
$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = 0.123995512547182
$ws.Range("AB14").Value = 551.6727934377117
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 551.6727934377117

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.00169565410619625
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"
